$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "NPCPoolsData"
$ws.Range("B10").Value = "NPCPoolsData.xlsx"
